$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4921406805515289
$ws.Range("B1").Value = 0.8842782378196716
$ws.Range("C1").Value = 5.871946811676025
$ws.Range("D1").Value = 1.583678722381592
$ws.Range("E1").Value = 1.17596447467804
